# Update the recalculated "continental/group aggregate" summary rows in Tab30.
# These rows (regional subtotal bands that follow every block of ~10 countries,
# plus the final "Etats fragiles" / fragile-states rows) were refreshed from the
# latest source computation; only the cached values move, no formulas exist on
# this values-only report sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab30")

# Row 13
$ws.Range("G13").Value = 3.82322159623598
$ws.Range("H13").Value = 6.0674919746544
$ws.Range("M13").Value = 2.31456772197732
$ws.Range("N13").Value = 4.41695357118348

# Row 23
$ws.Range("G23").Value = 3.10993318547215
$ws.Range("H23").Value = 6.85487707602731
$ws.Range("M23").Value = 5.7101125237588
$ws.Range("N23").Value = 13.0979600927452

# Row 38
$ws.Range("G38").Value = 1.3245840392668
$ws.Range("H38").Value = 3.30735168573735
$ws.Range("M38").Value = 6.16285483341604
$ws.Range("N38").Value = 17.0856424956151

# Row 45
$ws.Range("G45").Value = 8.42054149419929
$ws.Range("H45").Value = 21.0530827814597
$ws.Range("M45").Value = 0.11268289160202
$ws.Range("N45").Value = 0.27801480739247

# Row 61
$ws.Range("G61").Value = 3.79549628986822
$ws.Range("H61").Value = 7.66226294672631
$ws.Range("M61").Value = 7.78663496235225
$ws.Range("N61").Value = 15.7644264706401

# Row 62
$ws.Range("G62").Value = 5.12542088758775
$ws.Range("H62").Value = 11.6180583011232
$ws.Range("M62").Value = 3.58711145779411
$ws.Range("N62").Value = 7.92286675586102

# Row 63
$ws.Range("G63").Value = 4.95406803531666
$ws.Range("H63").Value = 9.98849669360542
$ws.Range("M63").Value = 1.35703724377071
$ws.Range("N63").Value = 2.67966789591237

# Row 64
$ws.Range("G64").Value = 2.83736886704326
$ws.Range("H64").Value = 7.54062078276273
$ws.Range("M64").Value = 0.76532731318104
$ws.Range("N64").Value = 2.05777550109396

# Row 65
$ws.Range("G65").Value = 8.56583519161135
$ws.Range("H65").Value = 16.7430854074291
$ws.Range("M65").Value = 3.40961334463872
$ws.Range("N65").Value = 6.72191085787484

# Row 66
$ws.Range("G66").Value = 4.96263367638807
$ws.Range("H66").Value = 10.0699557393625
$ws.Range("M66").Value = 1.46851490537229
$ws.Range("N66").Value = 2.94176659403806

# Row 67
$ws.Range("G67").Value = 6.32511014489715
$ws.Range("H67").Value = 15.809428496655
$ws.Range("J67").Value = 480.134942700188
$ws.Range("M67").Value = 2.28464538664843
$ws.Range("N67").Value = 6.40721088780889

# Row 68
$ws.Range("G68").Value = 6.62232320160944
$ws.Range("H68").Value = 15.3510155043476
$ws.Range("M68").Value = 3.80611525386085
$ws.Range("N68").Value = 7.66649338433339

# Row 69
$ws.Range("G69").Value = 1.29469953981599
$ws.Range("H69").Value = 3.35300049733527
$ws.Range("M69").Value = 6.4807802880922
$ws.Range("N69").Value = 16.7378705761721

# Row 70
$ws.Range("G70").Value = 2.63357050043448
$ws.Range("H70").Value = 6.1287397284837
$ws.Range("M70").Value = 4.79126728624281
$ws.Range("N70").Value = 11.5410540653308
$ws.Range("Z70").Value = 21.4169956665641

# Row 71
$ws.Range("G71").Value = 3.79549628986822
$ws.Range("H71").Value = 7.66226294672631
$ws.Range("M71").Value = 7.78663496235225
$ws.Range("N71").Value = 15.7644264706401

# Row 72
$ws.Range("G72").Value = 1.13128179482923
$ws.Range("H72").Value = 3.15235559461511
$ws.Range("M72").Value = 5.71028961611953
$ws.Range("N72").Value = 17.0970938361363

# Row 73
$ws.Range("G73").Value = 3.37162953502644
$ws.Range("H73").Value = 5.57555231070118
$ws.Range("J73").Value = 532.448906333688
$ws.Range("M73").Value = 3.32577597510224
$ws.Range("N73").Value = 7.27188064999034
$ws.Range("Z73").Value = 33.9032737136766

# Row 74
$ws.Range("G74").Value = 6.48072877276597
$ws.Range("H74").Value = 16.5562423959025
$ws.Range("M74").Value = 0.23114748739826
$ws.Range("N74").Value = 0.57049668741041

# Row 75
$ws.Range("G75").Value = 1.95533583111381
$ws.Range("H75").Value = 4.86124594514352
$ws.Range("M75").Value = 4.09671738435824
$ws.Range("N75").Value = 9.65978429079873

# Row 76
$ws.Range("G76").Value = 4.07650229126671
$ws.Range("H76").Value = 9.7396445491637
$ws.Range("M76").Value = 2.89381827226722
$ws.Range("N76").Value = 6.1258307243664

# Row 77
$ws.Range("G77").Value = 2.63850094615282
$ws.Range("H77").Value = 7.19821327998052
$ws.Range("M77").Value = 0.60565031069471
$ws.Range("N77").Value = 1.68465915585034

# Row 78
$ws.Range("G78").Value = 3.49513268572667
$ws.Range("H78").Value = 5.30140422599617
$ws.Range("M78").Value = 0.15682336127259
$ws.Range("N78").Value = 0.17434128409693

# Row 79
$ws.Range("G79").Value = 2.4403554600542
$ws.Range("H79").Value = 4.99586235115654
$ws.Range("M79").Value = 0.07676747876898
$ws.Range("N79").Value = 0.16870141105783

# Row 80
$ws.Range("G80").Value = 4.45210746503415
$ws.Range("H80").Value = 13.1094168915562
$ws.Range("M80").Value = 1.23703533698228
$ws.Range("N80").Value = 2.89760833179077

# Row 81
$ws.Range("G81").Value = 5.06102032321471
$ws.Range("H81").Value = 14.1379762638066
$ws.Range("M81").Value = 0.24674527579716
$ws.Range("N81").Value = 0.43664061701213

# Row 82
$ws.Range("G82").Value = 5.24768442924433
$ws.Range("H82").Value = 11.3472501202098
$ws.Range("M82").Value = 4.01384976430711
$ws.Range("N82").Value = 8.83537773824448

# Row 83
$ws.Range("G83").Value = 4.9441537424884
$ws.Range("H83").Value = 9.60384708626779
$ws.Range("M83").Value = 1.459959396365
$ws.Range("N83").Value = 2.8875926603268

# Row 84
$ws.Range("G84").Value = 1.23226436407572
$ws.Range("H84").Value = 3.24609621116279
$ws.Range("M84").Value = 7.52848753581852
$ws.Range("N84").Value = 19.9939699750505

# Row 85
$ws.Range("G85").Value = 3.00046066222706
$ws.Range("H85").Value = 9.0619583701789
$ws.Range("M85").Value = 5.62927433706551
$ws.Range("N85").Value = 17.0685796884566

# Row 86
$ws.Range("G86").Value = 6.02327974426446
$ws.Range("H86").Value = 14.0864428425844
$ws.Range("M86").Value = 3.47259706026037
$ws.Range("N86").Value = 7.12498888883807

# Row 87
$ws.Range("G87").Value = 6.70097153420191
$ws.Range("H87").Value = 13.6786456971891
$ws.Range("M87").Value = 4.53478467833385
$ws.Range("N87").Value = 9.28390553725479

# Row 88
$ws.Range("G88").Value = 4.71008685936634
$ws.Range("H88").Value = 8.58089437749012
$ws.Range("M88").Value = 0.84810433331845
$ws.Range("N88").Value = 1.31179200045941

# Row 89
$ws.Range("G89").Value = 7.76832789774173
$ws.Range("H89").Value = 14.7071989460874
$ws.Range("M89").Value = 1.72547920906354
$ws.Range("N89").Value = 3.27604493366546

# Row 90
$ws.Range("G90").Value = 2.42815987849561
$ws.Range("H90").Value = 5.49563103366654
$ws.Range("M90").Value = 0.05013520686908
$ws.Range("N90").Value = 0.06771991993603

# Row 91
$ws.Range("G91").Value = 1.54846044675996
$ws.Range("H91").Value = 3.9182305916538
$ws.Range("M91").Value = 6.99845399002857
$ws.Range("N91").Value = 17.7998118974033

# Row 92
$ws.Range("G92").Value = 5.17944822611925
$ws.Range("H92").Value = 11.6222496459853
$ws.Range("M92").Value = 8.19241241241595
$ws.Range("N92").Value = 17.895532762768

# Row 93
$ws.Range("G93").Value = 4.72313544079465
$ws.Range("H93").Value = 8.08946775735685
$ws.Range("M93").Value = 1.65771329728609
$ws.Range("N93").Value = 3.79762711144162

# Row 94
$ws.Range("G94").Value = 2.66251703071973
$ws.Range("H94").Value = 10.0085988976299
$ws.Range("M94").Value = 1.35115920744665
$ws.Range("N94").Value = 2.45174923674919

# Row 95
$ws.Range("G95").Value = 1.46649738322176
$ws.Range("H95").Value = 3.31533542694901
$ws.Range("M95").Value = 7.21589195456024
$ws.Range("N95").Value = 18.3696638836571

# Row 96
$ws.Range("G96").Value = 6.45905942827212
$ws.Range("H96").Value = 12.3146899523903
$ws.Range("M96").Value = 2.41028786556841
$ws.Range("N96").Value = 5.22979610961454

# Row 97
$ws.Range("C97").Value = 151565
$ws.Range("D97").Value = 194.530153488713
$ws.Range("E97").Value = 4.98901489022855
$ws.Range("F97").Value = 99.856380264071
$ws.Range("G97").Value = 2.87246635452774
$ws.Range("H97").Value = 6.35048540930384
$ws.Range("I97").Value = 502570
$ws.Range("J97").Value = 645.036909832893
$ws.Range("K97").Value = 16.9204508708274
$ws.Range("L97").Value = 238.893708118063
$ws.Range("M97").Value = 6.73575271481522
$ws.Range("N97").Value = 14.8564527656435
$ws.Range("O97").Value = 6.28133521046258
$ws.Range("P97").Value = 52.9241637522322
$ws.Range("Q97").Value = 28.7806826235008
$ws.Range("R97").Value = 40.7260007280727
$ws.Range("S97").Value = 1.47343130805293
$ws.Range("T97").Value = 486.052644912582
$ws.Range("U97").Value = 342.309666388725
$ws.Range("V97").Value = 322.565943556166
$ws.Range("W97").Value = 15.2169003382191
$ws.Range("X97").Value = 4.18798463255811
$ws.Range("Y97").Value = 50.4059823703508
$ws.Range("Z97").Value = 23.722500755002
$ws.Range("AA97").Value = 12.1044692214697

# Row 98
$ws.Range("C98").Value = 309098
$ws.Range("D98").Value = 463.218566831445
$ws.Range("E98").Value = 12.3684217244649
$ws.Range("F98").Value = 339.285232091155
$ws.Range("G98").Value = 5.45825549930173
$ws.Range("H98").Value = 13.1816843695761
$ws.Range("I98").Value = 339093
$ws.Range("J98").Value = 508.169491496467
$ws.Range("K98").Value = 12.6135478662757
$ws.Range("L98").Value = 289.166594807331
$ws.Range("M98").Value = 4.66554895556928
$ws.Range("N98").Value = 9.77902784409221
$ws.Range("O98").Value = 43.8659281836933
$ws.Range("P98").Value = 61.4776390242783
$ws.Range("Q98").Value = 129.576154366945
$ws.Range("R98").Value = 130.903926304042
$ws.Range("S98").Value = 4.36096652026058
$ws.Range("T98").Value = 163.844958978216
$ws.Range("U98").Value = 89.0731158298034
$ws.Range("V98").Value = 72.2331911603299
$ws.Range("W98").Value = 46.8406768238023
$ws.Range("X98").Value = 12.6408084530577
$ws.Range("Y98").Value = 155.365801393462
$ws.Range("Z98").Value = 54.7383722058549
$ws.Range("AA98").Value = 27.3586958054561
